$wb = $excel.ActiveWorkbook

# --- ProductLoanInput sheet: rename "Currency" row label/value and drop the
#     now-unused column C (empty, style-only cells) ---
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws1.Range("A6").Value = "currency"
$ws1.Range("B6").Value = "US Dollar"
$ws1.Columns("C:C").Delete()

# Update the selection on the input sheet to match the new active cell
$ws1.Range("A6:B6").Select()

# --- Make ProductLoanInput the active sheet/tab (was ProductLoanOutput) ---
$ws1.Activate()
